# Generate Report for Archive
# - Updates the localization "Status" text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F3 and the
#   per-locale "Status" column C2:C3 on the zh-cn / de-de sheets).
# - The status column(s) narrow accordingly (their width follows the new,
#   shorter text) on all three sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn.Range("C2:C3").Value = $newStatus
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe.Range("C2:C3").Value = $newStatus
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5

Write-Host "Updated status text '$oldStatus' -> '$newStatus' and resized status columns."
